# Add two new columns, I ("I0") and J ("IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - same bold/bordered style (s="1") as the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for column I (I0) and column J (IF), rows 2-70.
$iValues = @(5,7,7,8,8,9,7,7,6,9,7,8,8,4,6,5,6,6,7,8,9,6,8,6,8,7,9,11,8,6,7,7,6,8,4,8,9,7,9,9,3,9,8,8,4,10,2,7,5,9,5,6,4,4,7,9,7,8,9,7,7,8,8,8,5,6,9,7,2)
$jValues = @(5,7,7,8,8,9,7,7,6,9,7,8,9,5,7,6,6,6,7,9,9,6,8,6,8,7,9,11,9,6,8,7,7,8,5,8,9,8,9,9,4,9,8,8,4,12,3,7,6,9,6,7,5,5,7,9,8,8,9,7,7,8,8,8,5,6,9,7,2)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
